$d = $word.ActiveDocument

# Replace the title text first (while font size is still 16pt, find by old text)
$d.Content.Find.Execute(
    "План лечения для врача-стоматолога", $true, $false, $false, $false, $false,
    $true, 1, $false, "План лечения", 2
)

# Update the first paragraph's formatting: center alignment and font size 20pt (sz=40, szCs=40)
$p1 = $d.Paragraphs(1)
$p1.Range.ParagraphFormat.Alignment = 1
$p1.Range.Font.Size = 20
$p1.Range.Font.SizeBi = 20
